$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Update header row text (same cell positions, just text corrections)
$ws.Range("C1").Value = "veluG50"
$ws.Range("D1").Value = "veluG100"
$ws.Range("A1").Value = "Molécules"

# Reset selection to the default single-cell A1 (the commit removed the
# previous full-column B1:B1048576 selection, leaving the view on A1)
$ws.Range("A1").Select()
